$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain-text values such as "26.016.97" or "216.36".
# Assigning these strings directly would make Excel auto-convert the numeric-looking
# ones (e.g. "216.36") into real numbers, corrupting the intended text content.
# Force the whole Price column to a text format first so every assignment below is
# kept as text, then restore the default "Normal" style so no extra formatting lingers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "26.016.97"
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("D3").Value = "1.664.91"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "216.36"
$ws.Range("E5").Value = "  -0.56%  "
$ws.Range("D6").Value = "0.5096"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").Value = "21.63"
$ws.Range("E10").Value = "  -1.46%  "
$ws.Range("D11").Value = "0.07422"
$ws.Range("E11").Value = "  +1.82%  "
$ws.Range("D12").Value = "1.673.09"
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "4.497"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "0.000008521"
$ws.Range("E15").Value = "  +3.90%  "
$ws.Range("D16").Value = "64.21"
$ws.Range("E16").Value = "  -0.96%  "
$ws.Range("D17").Value = "26.108.84"
$ws.Range("E17").Value = "  -1.78%  "
$ws.Range("D18").Value = "4.888"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "10.73"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").Value = "188.63"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("E23").Value = "  +0.21%  "
$ws.Range("D24").Value = "145.61"
$ws.Range("E24").Value = "  +0.62%  "
$ws.Range("D25").Value = "7.602"
$ws.Range("E25").Value = "  +1.09%  "
$ws.Range("E26").Value = "  +4.40%  "
$ws.Range("D27").Value = "15.59"
$ws.Range("E27").Value = "  +0.67%  "
$ws.Range("D28").Value = "0.06472"
$ws.Range("E28").Value = "  +13.51%  "
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "3.522"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "3.502"
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("D33").Value = "1.625"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").Value = "0.6048"
$ws.Range("E35").Value = "  +1.94%  "
$ws.Range("D36").Value = "2.367"
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("D37").Value = "2.687"
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").Value = "6.190"
$ws.Range("E38").Value = "  +4.82%  "
$ws.Range("D39").Value = "0.01609"
$ws.Range("E39").Value = "  +0.94%  "
$ws.Range("D40").Value = "1.074.60"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Value = "0.8586"
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Value = "100.53"
$ws.Range("E43").Value = "  +2.51%  "
$ws.Range("D44").Value = "1.814.07"
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("E45").Value = "  +8.17%  "
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("D47").Value = "1.006"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").Value = "8.033"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("E49").Value = "  +0.17%  "
$ws.Range("D50").Value = "0.4294"
$ws.Range("E50").Value = "  -0.46%  "
$ws.Range("D51").Value = "5.933"
$ws.Range("E51").Value = "  +4.37%  "

$priceRange.Style = "Normal"
